$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (Text number format) on cells whose new values look like plain numbers,
# so Excel keeps them as text strings instead of auto-converting to numeric values
# (mirrors the original workbook, where these Price cells are stored as inline text).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.185.27'
$ws.Range("E2").Value = '  +1.59%  '
$ws.Range("D3").Value = '1.814.09'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '312.15'
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("D7").Value = '0.4608'
$ws.Range("E7").Value = '  +4.68%  '
$ws.Range("D8").Value = '0.3747'
$ws.Range("E8").Value = '  +2.29%  '
$ws.Range("D9").Value = '0.07390'
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("D10").Value = '0.8649'
$ws.Range("E10").Value = '  +1.64%  '
$ws.Range("D11").Value = '20.60'
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("D12").Value = '1.814.60'
$ws.Range("E12").Value = '  -14.27%  '
$ws.Range("D13").Value = '6.656'
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("D14").Value = '5.386'
$ws.Range("E14").Value = '  +2.51%  '
$ws.Range("D15").Value = '0.07082'
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").Value = '91.80'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  +0.23%  '
$ws.Range("D18").Value = '0.000008737'
$ws.Range("E18").Value = '  +1.31%  '
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").Value = '14.88'
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("D21").Value = '27.185.71'
$ws.Range("E21").Value = '  +1.48%  '
$ws.Range("D22").Value = '5.309'
$ws.Range("E22").Value = '  +3.48%  '
$ws.Range("D23").Value = '10.91'
$ws.Range("E23").Value = '  +1.18%  '
$ws.Range("D24").Value = '2.049.04'
$ws.Range("E24").Value = '  -5.22%  '
$ws.Range("E25").Value = '  -2.46%  '
$ws.Range("D26").Value = '151.86'
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("D27").Value = '2.222'
$ws.Range("E27").Value = '  +1.95%  '
$ws.Range("E28").Value = '  +1.11%  '
$ws.Range("D29").Value = '5.263'
$ws.Range("E29").Value = '  +1.96%  '
$ws.Range("D30").Value = '116.94'
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").Value = '0.08889'
$ws.Range("E31").Value = '  +1.35%  '
$ws.Range("D32").Value = '0.7720'
$ws.Range("E32").Value = '  +5.10%  '
$ws.Range("E33").Value = '  +1.88%  '
$ws.Range("D34").Value = '4.514'
$ws.Range("E34").Value = '  +2.35%  '
$ws.Range("D35").Value = '2.919'
$ws.Range("E35").Value = '  +0.65%  '
$ws.Range("D36").Value = '1.001'
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("D37").Value = '1.119'
$ws.Range("E37").Value = '  +3.63%  '
$ws.Range("D38").Value = '0.01960'
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("D39").Value = '0.05236'
$ws.Range("E39").Value = '  +1.40%  '
$ws.Range("D40").Value = '7.261'
$ws.Range("E40").Value = '  +3.51%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '2.380'
$ws.Range("E41").Value = '  +21.55%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").Value = '2.912'
$ws.Range("E42").Value = '  +3.97%  '
$ws.Range("D43").Value = '0.5281'
$ws.Range("E43").Value = '  +1.91%  '
$ws.Range("D44").Value = '0.1680'
$ws.Range("E44").Value = '  +0.72%  '
$ws.Range("D45").Value = '8.597'
$ws.Range("E45").Value = '  +2.53%  '
$ws.Range("D46").Value = '0.5020'
$ws.Range("E46").Value = '  +0.92%  '
$ws.Range("E47").Value = '  +1.47%  '
$ws.Range("D48").Value = '105.13'
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("D50").Value = '1.669'
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("D51").Value = '0.06318'
$ws.Range("E51").Value = '  +0.18%  '
